# Update cryptos list: prices (col D) and 1h volume/change (col E) for rows 2-51.
# Rows 43, 44, 45, 47, 48 also have their Coin name / Link swapped (re-ranking).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new D-column prices look like plain numbers (e.g. "231.91") and would
# otherwise be auto-converted to numeric values by the Value setter, which can silently
# drop significant trailing zeros (e.g. "0.760" -> 0.76). Force those cells to remain
# plain text first so the original formatting of the price string is preserved exactly.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '37.555.65'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '2.069.13'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '231.91'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = '0.628'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '57.43'
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').Value = '0.388'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').Value = '14.84'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '2.373.94'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').Value = '20.84'
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('D15').Value = '0.760'
$ws.Range('E15').Value = '  -2.19%  '
$ws.Range('D16').Value = '5.30'
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').Value = '2.069.44'
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').Value = '37.521.25'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').Value = '70.32'
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('D20').Value = '5.92'
$ws.Range('E20').Value = '  -3.46%  '
$ws.Range('D21').Value = '0.0₃0826'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('D22').Value = '227.35'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '2.34'
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('D26').Value = '9.59'
$ws.Range('E26').Value = '  +4.98%  '
$ws.Range('D27').Value = '168.63'
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('D28').Value = '0.132'
$ws.Range('E28').Value = '  -3.50%  '
$ws.Range('D29').Value = '19.38'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('E30').Value = '  -2.99%  '
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').Value = '4.57'
$ws.Range('E32').Value = '  -2.77%  '
$ws.Range('D33').Value = '0.0629'
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('D34').Value = '4.58'
$ws.Range('E34').Value = '  -2.44%  '
$ws.Range('D35').Value = '2.46'
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('E37').Value = '  -3.70%  '
$ws.Range('D38').Value = '0.995'
$ws.Range('E38').Value = '  -0.63%  '
$ws.Range('D39').Value = '5.27'
$ws.Range('E39').Value = '  -1.87%  '
$ws.Range('D40').Value = '0.0229'
$ws.Range('E40').Value = '  +6.02%  '
$ws.Range('D41').Value = '99.07'
$ws.Range('E41').Value = '  -0.87%  '
$ws.Range('D42').Value = '0.0959'
$ws.Range('E42').Value = '  -1.87%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.20'
$ws.Range('E43').Value = '  +4.10%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.479.71'
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').Value = '2.88'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '16.55'
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '1.03'
$ws.Range('E47').Value = '  -2.89%  '
$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').Value = '3.99'
$ws.Range('E48').Value = '  -4.69%  '
$ws.Range('D49').Value = '7.20'
$ws.Range('E49').Value = '  -2.99%  '
$ws.Range('D50').Value = '2.95'
$ws.Range('E50').Value = '  -1.96%  '
$ws.Range('D51').Value = '2.256.43'
$ws.Range('E51').Value = '  -0.82%  '
